# Updated cryptos list - refresh Price / Volume(1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.334.82"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.393.13"
$ws.Range("E3").Value = "  +7.45%  "
$ws.Range("E4").Value = "  -0.19%  "
$origStyleD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.84"
$ws.Range("D5").Style = $origStyleD5
$ws.Range("E5").Value = "  +10.45%  "
$origStyleD6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.73"
$ws.Range("D6").Style = $origStyleD6
$ws.Range("E6").Value = "  -5.50%  "
$ws.Range("E7").Value = "  +4.43%  "
$ws.Range("E8").Value = "  -0.08%  "
$origStyleD9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.652"
$ws.Range("D9").Style = $origStyleD9
$ws.Range("E9").Value = "  +8.67%  "
$origStyleD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.97"
$ws.Range("D10").Style = $origStyleD10
$ws.Range("E10").Value = "  -4.27%  "
$origStyleD11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("D11").Style = $origStyleD11
$ws.Range("E11").Value = "  +3.29%  "
$origStyleD12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.60"
$ws.Range("D12").Style = $origStyleD12
$ws.Range("E12").Value = "  -0.68%  "
$origStyleD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.55"
$ws.Range("D13").Style = $origStyleD13
$ws.Range("E13").Value = "  +17.40%  "
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "2.753.82"
$ws.Range("E16").Value = "  +7.47%  "
$ws.Range("D17").Value = "2.394.61"
$ws.Range("E17").Value = "  +7.39%  "
$ws.Range("D18").Value = "43.327.36"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("E20").Value = "  +3.93%  "
$origStyleD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.43"
$ws.Range("D21").Style = $origStyleD21
$ws.Range("E21").Value = "  +3.91%  "
$origStyleD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "271.86"
$ws.Range("D22").Style = $origStyleD22
$ws.Range("E22").Value = "  +15.40%  "
$origStyleD23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.43"
$ws.Range("D23").Style = $origStyleD23
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +11.52%  "
$origStyleD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.84"
$ws.Range("D26").Style = $origStyleD26
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("E27").Value = "  +0.07%  "
$origStyleD28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.92"
$ws.Range("D28").Style = $origStyleD28
$ws.Range("E28").Value = "  +7.79%  "
$origStyleD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "176.92"
$ws.Range("D29").Style = $origStyleD29
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("E30").Value = "  -0.37%  "
$origStyleD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.75"
$ws.Range("D31").Style = $origStyleD31
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +2.31%  "
$origStyleD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0931"
$ws.Range("D33").Style = $origStyleD33
$ws.Range("E33").Value = "  +6.00%  "
$origStyleD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.91"
$ws.Range("D34").Style = $origStyleD34
$ws.Range("E34").Value = "  +5.13%  "
$ws.Range("E35").Value = "  +6.42%  "
$origStyleD36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.86"
$ws.Range("D36").Style = $origStyleD36
$ws.Range("E36").Value = "  -3.07%  "
$origStyleD37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.13"
$ws.Range("D37").Style = $origStyleD37
$ws.Range("E37").Value = "  -1.93%  "
$origStyleD38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0368"
$ws.Range("D38").Style = $origStyleD38
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("E39").Value = "  +4.30%  "
$origStyleD40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("D40").Style = $origStyleD40
$ws.Range("E40").Value = "  +18.36%  "
$ws.Range("E41").Value = "  +21.94%  "
$origStyleD42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "127.36"
$ws.Range("D42").Style = $origStyleD42
$ws.Range("E42").Value = "  +26.10%  "
$ws.Range("E43").Value = "  +1.57%  "
$origStyleD44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.43"
$ws.Range("D44").Style = $origStyleD44
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +2.23%  "
$origStyleD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.65"
$ws.Range("D47").Style = $origStyleD47
$ws.Range("E47").Value = "  +14.54%  "
$origStyleD48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.01"
$ws.Range("D48").Style = $origStyleD48
$ws.Range("E48").Value = "  +61.78%  "
$origStyleD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.65"
$ws.Range("D49").Style = $origStyleD49
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").Value = "1.606.72"
$ws.Range("E51").Value = "  +12.46%  "
